$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue List")

# Mark issues as resolved (待解决 -> 已解决), with a green "resolved" fill.
$resolvedRows = @(3, 4, 5, 6, 7, 9)
foreach ($r in $resolvedRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = "已解决"
    $cell.Interior.Color = 5296274
    $cell.VerticalAlignment = -4160
}

# Row 8 (维护页面的备注输入框...) is put on hold, with a blue fill.
$holdCell = $ws.Cells.Item(8, 4)
$holdCell.Value = "Hold"
$holdCell.Interior.Color = 15773696
$holdCell.VerticalAlignment = -4160

# Restore the active selection to D8, matching the author's last edit location.
$ws.Range("D8").Select()
